# Updated cryptos list on Fri Mar 15 04:45:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.486.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.51%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.671.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -7.52%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.57%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.662.74"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -7.55%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -7.75%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - USDC
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.709"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -11.74%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -8.26%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -12.19%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.262.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.62%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.718.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.98%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - TRON
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -8.10%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - Polygon
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.11%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - WrappedBTC
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.508.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.44%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.19%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - PancakeSwap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.44%  "
$ws.Range("E23").Style = "Normal"

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.79%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - ImmutableX
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -9.16%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -10.20%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("E27").Style = "Normal"

# Row 28 - Toncoin
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.20%  "
$ws.Range("E28").Style = "Normal"

# Row 29 - LEO
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -9.05%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.13%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - Cosmos
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.31%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "OKB"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.72%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Hedera"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.35%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "44.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.29%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - PEPE
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0915"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -10.37%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - Bittensor
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "595.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.65%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - TheGraph
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - Dai
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.59%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - Kaspa
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.135"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.18%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - ThetaToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -12.53%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - VeChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0439"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.28%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - THORChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -12.20%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - Fetch.AI
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - Stellar
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -9.48%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - Maker
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.736.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("E49").Style = "Normal"

# Row 50 - WEMIXToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -17.87%  "
$ws.Range("E50").Style = "Normal"

# Row 51 - ApeXProtocol
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.29%  "
$ws.Range("E51").Style = "Normal"

